# The workbook contains a daily price log for "Acelga" at the Chillán
# terminal, ordered from most-recent to oldest date. A new daily record
# (fecha 2022-02-09, serial 44601) was inserted right above the existing
# row for serial 44589, pushing every subsequent record down by one row.
#
# Concretely: insert a new row at row 103 (which shifts the old rows
# 103..207 down to 104..208, growing the sheet from 207 to 208 data rows)
# and populate the newly inserted row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 103 (and everything below it) down by one row.
$ws.Rows("103:103").Insert()

# Populate the newly inserted row 103 with the new daily record.
$ws.Range("A103").Value = 7
$ws.Range("B103").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C103").Value = "Ñuble"
$ws.Range("D103").Value = 44601
$ws.Range("E103").Value = 16
$ws.Range("F103").Value = 100112009
$ws.Range("G103").Value = "Acelga"
$ws.Range("H103").Value = "Sin especificar"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 100
$ws.Range("K103").Value = 400
$ws.Range("L103").Value = 450
$ws.Range("M103").Value = 425
$ws.Range("N103").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O103").Value = "Provincia de Diguillín"
$ws.Range("P103").Value = 425
$ws.Range("Q103").Value = 1
$ws.Range("R103").Value = "Hortaliza"
